$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 10818.667
$ws.Range("I38").Value = 13425.223
$ws.Range("K38").Value = 40275.669
$ws.Range("M38").Value = -39903.669
$ws.Range("H43").Value = 83335170
$ws.Range("J43").Value = 125002250
$ws.Range("L43").Value = 125002250
$ws.Range("N43").Value = -125002388
$ws.Range("H53").Value = 534
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 534
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 534
$ws.Range("M53").ClearContents()
$ws.Range("N53").Value = -1808
$ws.Range("H86").Value = 2192.7144
$ws.Range("I86").Value = 3100
$ws.Range("K86").Value = 3100
$ws.Range("M86").Value = -1977
$ws.Range("H89").Value = 2192.7144
$ws.Range("I89").Value = 3100
$ws.Range("K89").Value = 15500
$ws.Range("M89").Value = -9884
$ws.Range("H92").Value = 1583.409
$ws.Range("I92").Value = 1490.8889
$ws.Range("J92").Value = 1999.75
$ws.Range("K92").Value = 1490.8889
$ws.Range("L92").Value = 1999.75
$ws.Range("M92").Value = -242.8888999999999
$ws.Range("N92").Value = -4495.75
$ws.Range("H96").Value = 5367.9
$ws.Range("I96").Value = 242.71428
$ws.Range("K96").Value = 728.14284
$ws.Range("M96").Value = 644.85716
$ws.Range("H98").Value = 5676.0757
$ws.Range("I98").Value = 5098.4897
$ws.Range("K98").Value = 5098.4897
$ws.Range("M98").Value = -3600.4897
$ws.Range("H111").Value = 990.25
$ws.Range("I111").Value = 1349
$ws.Range("J111").Value = 631.5
$ws.Range("K111").Value = 4047
$ws.Range("L111").Value = 1894.5
$ws.Range("M111").Value = -980
$ws.Range("N111").Value = -8028.5
$ws.Range("H122").Value = 5676.0757
$ws.Range("I122").Value = 5098.4897
$ws.Range("K122").Value = 15295.4691
$ws.Range("M122").Value = -12845.4691
$ws.Range("H132").Value = 1696722.8
$ws.Range("I132").Value = 1721.196
$ws.Range("K132").Value = 5163.588
$ws.Range("M132").Value = -2633.588
$ws.Range("H137").Value = 1088169.2
$ws.Range("I137").Value = 1018.7143
$ws.Range("J137").Value = 3624853.8
$ws.Range("K137").Value = 3056.1429
$ws.Range("L137").Value = 10874561.4
$ws.Range("M137").Value = -506.1428999999998
$ws.Range("N137").Value = -10879661.4

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1003.2
$ws.Range("J2").Value = 2642
$ws.Range("L2").Value = 2642
$ws.Range("N2").Value = -2868
$ws.Range("H32").Value = 123374.336
$ws.Range("I32").Value = 123374.336
$ws.Range("K32").Value = 123374.336
$ws.Range("M32").Value = -123087.336
$ws.Range("H45").Value = 2654.8333
$ws.Range("J45").Value = 3751.6667
$ws.Range("L45").Value = 3751.6667
$ws.Range("N45").Value = -4505.6667
$ws.Range("H74").Value = 1632.5846
$ws.Range("I74").Value = 1323.5536
$ws.Range("K74").Value = 1323.5536
$ws.Range("M74").Value = -449.5536
$ws.Range("H77").Value = 1632.5846
$ws.Range("I77").Value = 1323.5536
$ws.Range("K77").Value = 6617.768
$ws.Range("M77").Value = -2249.768
$ws.Range("H82").Value = 200047040
$ws.Range("J82").Value = 250050000
$ws.Range("L82").Value = 250050000
$ws.Range("N82").Value = -250050722
$ws.Range("H85").Value = 200047040
$ws.Range("J85").Value = 250050000
$ws.Range("L85").Value = 250050000
$ws.Range("N85").Value = -250052496
$ws.Range("H92").Value = 35700
$ws.Range("I92").Value = 29900
$ws.Range("J92").Value = 36666.668
$ws.Range("K92").Value = 29900
$ws.Range("L92").Value = 36666.668
$ws.Range("M92").Value = -27404
$ws.Range("N92").Value = -41658.668
$ws.Range("H94").Value = 20407.334
$ws.Range("J94").Value = 20407.334
$ws.Range("L94").Value = 20407.334
$ws.Range("N94").Value = -22209.334
$ws.Range("H95").Value = 77919
$ws.Range("J95").Value = 77919
$ws.Range("L95").Value = 77919
$ws.Range("N95").Value = -83411
$ws.Range("H101").Value = 28000
$ws.Range("J101").Value = 28000
$ws.Range("L101").Value = 28000
$ws.Range("N101").Value = -34490
$ws.Range("H102").Value = 3805.8667
$ws.Range("I102").Value = 1529.6875
$ws.Range("J102").Value = 6407.2144
$ws.Range("K102").Value = 1529.6875
$ws.Range("L102").Value = 6407.2144
$ws.Range("M102").Value = 92.3125
$ws.Range("N102").Value = -9651.214400000001
$ws.Range("H110").Value = 5977.0713
$ws.Range("I110").Value = 6052.231
$ws.Range("K110").Value = 6052.231
$ws.Range("M110").Value = -4007.231
$ws.Range("H116").Value = 1003.2
$ws.Range("J116").Value = 2642
$ws.Range("L116").Value = 2642
$ws.Range("N116").Value = -7230
$ws.Range("H132").Value = 3367.1365
$ws.Range("I132").Value = 3528.275
$ws.Range("K132").Value = 10584.825
$ws.Range("M132").Value = -8054.825000000001

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1003.2
$ws.Range("J3").Value = 2642
$ws.Range("L3").Value = 2642
$ws.Range("N3").Value = -2870
$ws.Range("H92").Value = 35446
$ws.Range("J92").Value = 35446
$ws.Range("L92").Value = 35446
$ws.Range("N92").Value = -40438
$ws.Range("H97").Value = 15999.75
$ws.Range("I97").Value = 12499.5
$ws.Range("J97").Value = 19500
$ws.Range("K97").Value = 12499.5
$ws.Range("L97").Value = 19500
$ws.Range("M97").Value = -11508.5
$ws.Range("N97").Value = -21482
$ws.Range("H101").Value = 37449.5
$ws.Range("J101").Value = 37449.5
$ws.Range("L101").Value = 37449.5
$ws.Range("N101").Value = -43939.5
$ws.Range("H102").Value = 9425.666999999999
$ws.Range("I102").Value = 8888.5
$ws.Range("J102").Value = 10500
$ws.Range("K102").Value = 8888.5
$ws.Range("L102").Value = 10500
$ws.Range("M102").Value = -5643.5
$ws.Range("N102").Value = -16990
$ws.Range("H103").Value = 42250
$ws.Range("J103").Value = 42250
$ws.Range("L103").Value = 42250
$ws.Range("N103").Value = -44594
$ws.Range("H105").Value = 2812.125
$ws.Range("I105").Value = 1949.6
$ws.Range("J105").Value = 4249.6665
$ws.Range("K105").Value = 1949.6
$ws.Range("L105").Value = 4249.6665
$ws.Range("M105").Value = -202.5999999999999
$ws.Range("N105").Value = -7743.6665

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3648.5
$ws.Range("I16").Value = 2849.5
$ws.Range("K16").Value = 2849.5
$ws.Range("M16").Value = -2562.5
$ws.Range("H19").Value = 566.2222
$ws.Range("I19").Value = 578
$ws.Range("K19").Value = 578
$ws.Range("M19").Value = -408
$ws.Range("H24").Value = 566.2222
$ws.Range("I24").Value = 578
$ws.Range("K24").Value = 578
$ws.Range("M24").Value = -408
$ws.Range("H29").Value = 3824.75
$ws.Range("J29").Value = 3824.75
$ws.Range("L29").Value = 3824.75
$ws.Range("N29").Value = -4410.75
$ws.Range("H42").Value = 6055.1665
$ws.Range("I42").Value = 4750
$ws.Range("J42").Value = 6707.75
$ws.Range("K42").Value = 4750
$ws.Range("L42").Value = 6707.75
$ws.Range("M42").Value = -4157
$ws.Range("N42").Value = -7893.75
$ws.Range("H43").Value = 16635.666
$ws.Range("J43").Value = 16635.666
$ws.Range("L43").Value = 16635.666
$ws.Range("N43").Value = -17003.666
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H58").Value = 2586.0625
$ws.Range("I58").Value = 2425.1333
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 2425.1333
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -2222.1333
$ws.Range("N58").Value = -5406
$ws.Range("H74").Value = 50210
$ws.Range("J74").Value = 50210
$ws.Range("L74").Value = 50210
$ws.Range("N74").Value = -51958
$ws.Range("H77").Value = 50210
$ws.Range("J77").Value = 50210
$ws.Range("L77").Value = 150630
$ws.Range("N77").Value = -159366
$ws.Range("H101").Value = 16635.666
$ws.Range("J101").Value = 16635.666
$ws.Range("L101").Value = 16635.666
$ws.Range("N101").Value = -23125.666
$ws.Range("H102").Value = 115241
$ws.Range("J102").Value = 115241
$ws.Range("L102").Value = 115241
$ws.Range("N102").Value = -120109
$ws.Range("H103").Value = 14666.667
$ws.Range("J103").Value = 14000
$ws.Range("L103").Value = 14000
$ws.Range("N103").Value = -16344
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H113").Value = 3648.5
$ws.Range("I113").Value = 2849.5
$ws.Range("K113").Value = 2849.5
$ws.Range("M113").Value = -679.5
$ws.Range("H132").Value = 1782.3823
$ws.Range("I132").Value = 1264.9231
$ws.Range("K132").Value = 3794.7693
$ws.Range("M132").Value = -1264.7693
$ws.Range("H136").Value = 2586.0625
$ws.Range("I136").Value = 2425.1333
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7275.3999
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -4725.3999
$ws.Range("N136").Value = -20100

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 50.133335
$ws.Range("I2").Value = 21.5
$ws.Range("J2").Value = 82.85714
$ws.Range("K2").Value = 129
$ws.Range("L2").Value = 497.14284
$ws.Range("M2").Value = -16
$ws.Range("N2").Value = -723.14284
$ws.Range("H5").Value = 1101.2632
$ws.Range("I5").Value = 751.875
$ws.Range("J5").Value = 2964.6667
$ws.Range("K5").Value = 2255.625
$ws.Range("L5").Value = 8894.000100000001
$ws.Range("M5").Value = -2143.625
$ws.Range("N5").Value = -9118.000100000001
$ws.Range("H33").Value = 279.6
$ws.Range("J33").Value = 316
$ws.Range("L33").Value = 1896
$ws.Range("N33").Value = -2462
$ws.Range("H44").Value = 466.4762
$ws.Range("I44").Value = 448.55554
$ws.Range("J44").Value = 479.91666
$ws.Range("K44").Value = 1345.66662
$ws.Range("L44").Value = 1439.74998
$ws.Range("M44").Value = -947.66662
$ws.Range("N44").Value = -2235.74998
$ws.Range("H49").Value = 2014.3334
$ws.Range("I49").Value = 1852.6666
$ws.Range("J49").Value = 2499.3333
$ws.Range("K49").Value = 5557.9998
$ws.Range("L49").Value = 7497.999899999999
$ws.Range("M49").Value = -5401.9998
$ws.Range("N49").Value = -7809.999899999999
$ws.Range("H68").Value = 2059.6667
$ws.Range("J68").Value = 2059.6667
$ws.Range("L68").Value = 6179.000100000001
$ws.Range("N68").Value = -7801.000100000001
$ws.Range("H71").Value = 2059.6667
$ws.Range("J71").Value = 2059.6667
$ws.Range("L71").Value = 18537.0003
$ws.Range("N71").Value = -26649.0003
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H92").Value = 365.66666
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H94").Value = 9151.875
$ws.Range("J94").Value = 10455.857
$ws.Range("L94").Value = 31367.571
$ws.Range("N94").Value = -32719.571
$ws.Range("H107").Value = 2279.4583
$ws.Range("J107").Value = 1857.2941
$ws.Range("L107").Value = 5571.8823
$ws.Range("N107").Value = -9411.882300000001
$ws.Range("H122").Value = 1025.9166
$ws.Range("J122").Value = 1117
$ws.Range("L122").Value = 10053
$ws.Range("N122").Value = -14953
$ws.Range("H130").Value = 10015
$ws.Range("I130").Value = 5030
$ws.Range("J130").Value = 15000
$ws.Range("K130").Value = 15090
$ws.Range("L130").Value = 45000
$ws.Range("M130").Value = -10070
$ws.Range("N130").Value = -55040
$ws.Range("H132").Value = 1622.6774
$ws.Range("I132").Value = 1082.15
$ws.Range("J132").Value = 2605.4546
$ws.Range("K132").Value = 9739.35
$ws.Range("L132").Value = 23449.0914
$ws.Range("M132").Value = -7209.35
$ws.Range("N132").Value = -28509.0914
$ws.Range("H135").Value = 1101.2632
$ws.Range("I135").Value = 751.875
$ws.Range("J135").Value = 2964.6667
$ws.Range("K135").Value = 6766.875
$ws.Range("L135").Value = 26682.0003
$ws.Range("M135").Value = -4231.875
$ws.Range("N135").Value = -31752.0003

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 30559.75
$ws.Range("J49").Value = 30559.75
$ws.Range("L49").Value = 30559.75
$ws.Range("N49").Value = -30927.75
$ws.Range("H70").Value = 17546626
$ws.Range("I70").Value = 20836038
$ws.Range("K70").Value = 20836038
$ws.Range("M70").Value = -20835768
$ws.Range("H73").Value = 17546626
$ws.Range("I73").Value = 20836038
$ws.Range("K73").Value = 20836038
$ws.Range("M73").Value = -20835102
$ws.Range("H102").Value = 1022.03705
$ws.Range("I102").Value = 735.24
$ws.Range("K102").Value = 735.24
$ws.Range("M102").Value = 886.76
$ws.Range("H122").Value = 8027
$ws.Range("I122").Value = 9732.076999999999
$ws.Range("K122").Value = 29196.231
$ws.Range("M122").Value = -26746.231
$ws.Range("H126").Value = 2350.923
$ws.Range("I126").Value = 1796.25
$ws.Range("J126").Value = 3238.4
$ws.Range("K126").Value = 5388.75
$ws.Range("L126").Value = 9715.200000000001
$ws.Range("M126").Value = -2918.75
$ws.Range("N126").Value = -14655.2
$ws.Range("H132").Value = 34484492
$ws.Range("I132").Value = 34484492
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 103453476
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -103450946
$ws.Range("N132").ClearContents()
$ws.Range("H133").Value = 94991.336
$ws.Range("J133").Value = 94991.336
$ws.Range("L133").Value = 94991.336
$ws.Range("N133").Value = -105111.336
$ws.Range("H134").Value = 73536.836
$ws.Range("J134").Value = 73536.836
$ws.Range("L134").Value = 220610.508
$ws.Range("N134").Value = -225680.508

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1557.9
$ws.Range("I55").Value = 1484.5714
$ws.Range("J55").Value = 1729
$ws.Range("K55").Value = 1484.5714
$ws.Range("L55").Value = 1729
$ws.Range("M55").Value = -1311.5714
$ws.Range("N55").Value = -2075
$ws.Range("H61").Value = 6406.923
$ws.Range("I61").Value = 4323.5
$ws.Range("J61").Value = 7332.8887
$ws.Range("K61").Value = 4323.5
$ws.Range("L61").Value = 7332.8887
$ws.Range("M61").Value = -4121.5
$ws.Range("N61").Value = -7736.8887
$ws.Range("H93").Value = 3338.151
$ws.Range("I93").Value = 1577.36
$ws.Range("J93").Value = 4910.2856
$ws.Range("K93").Value = 1577.36
$ws.Range("L93").Value = 4910.2856
$ws.Range("M93").Value = -329.3599999999999
$ws.Range("N93").Value = -7406.2856
$ws.Range("H113").Value = 6406.923
$ws.Range("I113").Value = 4323.5
$ws.Range("J113").Value = 7332.8887
$ws.Range("K113").Value = 4323.5
$ws.Range("L113").Value = 7332.8887
$ws.Range("M113").Value = -2153.5
$ws.Range("N113").Value = -11672.8887
$ws.Range("H132").Value = 2899.04
$ws.Range("I132").Value = 2216.9048
$ws.Range("K132").Value = 6650.714399999999
$ws.Range("M132").Value = -4120.714399999999
$ws.Range("H136").Value = 23901.773
$ws.Range("I136").Value = 3429.9546
$ws.Range("J136").Value = 73944
$ws.Range("K136").Value = 10289.8638
$ws.Range("L136").Value = 221832
$ws.Range("M136").Value = -7739.863799999999
$ws.Range("N136").Value = -226932

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 23110.889
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 25874.75
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 25874.75
$ws.Range("M54").Value = -480
$ws.Range("N54").Value = -26914.75
$ws.Range("H62").Value = 10428213
$ws.Range("J62").Value = 13825
$ws.Range("L62").Value = 13825
$ws.Range("N62").Value = -15073
$ws.Range("H65").Value = 10428213
$ws.Range("J65").Value = 13825
$ws.Range("L65").Value = 69125
$ws.Range("N65").Value = -75365
$ws.Range("H96").Value = 3691.4285
$ws.Range("J96").Value = 4333.727
$ws.Range("L96").Value = 4333.727
$ws.Range("N96").Value = -7079.727
$ws.Range("H107").Value = 26316252
$ws.Range("J107").Value = 55555996
$ws.Range("L107").Value = 166667988
$ws.Range("N107").Value = -166671828
$ws.Range("H122").Value = 455629.34
$ws.Range("I122").Value = 974425.4
$ws.Range("J122").Value = 6006.1333
$ws.Range("K122").Value = 2923276.2
$ws.Range("L122").Value = 18018.3999
$ws.Range("M122").Value = -2920826.2
$ws.Range("N122").Value = -22918.3999
